# Updates the crypto price/volume table (columns D=Price, E=Volume(1h))
# with the latest scraped figures, plus a swap of the ImmutableX /
# InternetComputer(DFINITY) rows (20-21).
#
# Note: some new Price values (column D) are plain decimals like "357.00"
# which Excel would normally auto-convert to a number, dropping the
# formatting (e.g. trailing zero). Prefixing with a leading apostrophe
# forces Excel to keep them as literal text, matching the original sheet's
# text-based Price column (the apostrophe itself is not stored in the
# cell's value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.722.17'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = '2.774.41'
$ws.Range('E3').Value = '  -1.04%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '''357.00'
$ws.Range('E5').Value = '  +1.44%  '

$ws.Range('D6').Value = '''109.06'
$ws.Range('E6').Value = '  -1.59%  '

$ws.Range('E7').Value = '  -0.72%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '''0.590'
$ws.Range('E9').Value = '  -0.60%  '

$ws.Range('D10').Value = '''39.75'
$ws.Range('E10').Value = '  -1.67%  '

$ws.Range('E11').Value = '  +2.59%  '

$ws.Range('D12').Value = '''0.0845'
$ws.Range('E12').Value = '  -0.50%  '

$ws.Range('D13').Value = '''19.46'
$ws.Range('E13').Value = '  -1.32%  '

$ws.Range('D14').Value = '''7.60'
$ws.Range('E14').Value = '  -1.68%  '

$ws.Range('D15').Value = '3.208.73'
$ws.Range('E15').Value = '  -1.18%  '

$ws.Range('D16').Value = '2.778.96'
$ws.Range('E16').Value = '  -1.33%  '

$ws.Range('D17').Value = '''0.934'
$ws.Range('E17').Value = '  +2.37%  '

$ws.Range('D18').Value = '51.634.91'
$ws.Range('E18').Value = '  +0.16%  '

$ws.Range('D19').Value = '''7.45'
$ws.Range('E19').Value = '  -0.69%  '

$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = '''13.12'
$ws.Range('E20').Value = '  -1.03%  '

$ws.Range('B21').Value = 'ImmutableX'
$ws.Range('C21').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D21').Value = '''3.01'
$ws.Range('E21').Value = '  -2.60%  '

$ws.Range('D22').Value = '0.0₃0970'
$ws.Range('E22').Value = '  -1.59%  '

$ws.Range('D23').Value = '''70.19'
$ws.Range('E23').Value = '  +0.57%  '

$ws.Range('D24').Value = '''268.64'
$ws.Range('E24').Value = '  +0.85%  '

$ws.Range('E25').Value = '  -1.52%  '

$ws.Range('D26').Value = '''26.37'
$ws.Range('E26').Value = '  -1.49%  '

$ws.Range('E27').Value = '  -0.01%  '

$ws.Range('D28').Value = '''0.163'
$ws.Range('E28').Value = '  +15.92%  '

$ws.Range('E29').Value = '  +0.24%  '

$ws.Range('E30').Value = '  -1.30%  '

$ws.Range('D31').Value = '''6.21'
$ws.Range('E31').Value = '  +6.35%  '

$ws.Range('E32').Value = '  -0.42%  '

$ws.Range('D33').Value = '''34.73'
$ws.Range('E33').Value = '  +1.91%  '

$ws.Range('D34').Value = '''0.0452'
$ws.Range('E34').Value = '  -8.48%  '

$ws.Range('D35').Value = '''0.0839'
$ws.Range('E35').Value = '  -0.22%  '

$ws.Range('D36').Value = '''5.14'
$ws.Range('E36').Value = '  -5.16%  '

$ws.Range('E37').Value = '  +0.14%  '

$ws.Range('D38').Value = '''18.72'
$ws.Range('E38').Value = '  +3.23%  '

$ws.Range('D39').Value = '''3.13'
$ws.Range('E39').Value = '  -3.44%  '

$ws.Range('E40').Value = '  -2.92%  '

$ws.Range('E41').Value = '  +2.24%  '

$ws.Range('E42').Value = '  -1.98%  '

$ws.Range('E43').Value = '  -1.98%  '

$ws.Range('D44').Value = '''119.61'
$ws.Range('E44').Value = '  -4.87%  '

$ws.Range('D45').Value = '''21.62'
$ws.Range('E45').Value = '  -6.56%  '

$ws.Range('D46').Value = '2.084.58'
$ws.Range('E46').Value = '  -0.02%  '

$ws.Range('E47').Value = '  -0.86%  '

$ws.Range('E48').Value = '  +1.01%  '

$ws.Range('D49').Value = '''5.74'
$ws.Range('E49').Value = '  -3.17%  '

$ws.Range('D50').Value = '''0.940'
$ws.Range('E50').Value = '  -3.54%  '

$ws.Range('E51').Value = '  +2.11%  '
